# Re #810 fixed combine_pow and crash in combining files together
#
# Swap the order of slide 5 ("SQW and DND objects. Composition:") and
# slide 6 ("Main Horace' objects interfaces:") so that the "Main
# Horace" slide now comes before the "SQW and DND" slide.

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(5)
$slide.MoveTo(6)
